$d = $word.ActiveDocument
$zwsp = [char]0x200B

# --- Paragraph ("-" + "Aportar en el desarrollo...") ---
# Both original runs get collapsed into a single new run with new wording.
$old1 = "-Aportar en el desarrollo de un sistema de inventarios y creación de un sitio WEB adecuado, favoreciendo a sus integrantes y principalmente al gerente de la empresa para contar con una interfaz adecuada de información que permita tener el control de los bienes con los que la empresa cuenta apoyado de un programa eficaz y verídico que brinden solución a las falencias existentes de la compañía y contar con una base de datos y/o programa que mantenga dicha información sin riesgo de que esta se extravié.  "
$new1 = "El alcance del proyecto será la creación e implementación de un sistema de inventario para la empresa `"SanBra Seguridad y Seguros`". Con la construcción de este sistema se espera lograr la organización de inventario y hacer un seguimiento en tiempo real de los activos, pasivos y patrimonios de la empresa, relacionados con la distribución de insumos de seguridad industrial, seguros de vida y seguros de vehículos que ofrece la empresa."

$range1 = $d.Content
$found1 = $range1.Find.Execute($old1)
if ($found1) {
    $range1.Text = $new1
}

# --- Paragraph ("-" + "Permitir que la empresa...") ---
# Locate this paragraph by searching for its leading run text, so the
# script does not depend on a hard-coded paragraph index.
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "-Permitir que la empresa cuente*") {
        $targetParaIndex = $i
    }
}

if ($targetParaIndex -gt 0) {
    $p5 = $d.Paragraphs($targetParaIndex)

    # The leading "-" run becomes the start of a new sentence ("El tiempo
    # prolongado ... para la empre"), and the second run's text is fully
    # replaced with "sa." to finish the sentence. These must stay as two
    # distinct runs (matching formatting), so toggle a formatting property
    # on/off after each edit to stop the engine from silently re-merging
    # the adjacent, identically-formatted runs back together.
    $dashRange = $d.Range($p5.Range.Start, $p5.Range.Start + 1)
    if ($dashRange.Text -eq "-") {
        $dashRange.Text = "El tiempo prolongado en el que lleváremos a cabo el proyecto será de 6 meses para poder mostrar avances frente al sistema de inventario que tenemos planteado para la empre"
        $dashRange.Font.Bold = 1
        $dashRange.Font.Bold = 0
    }

    $old2b = "Permitir que la empresa cuente con un sistema adecuado favorecerá la economía de la misma, ya que se tendrá claro la cantidad de insumos existentes y los faltantes sin que ocurra el gasto innecesario de productos, maquinarias, equipos, herramientas, entre otros. " + $zwsp
    $range2b = $d.Content
    $found2b = $range2b.Find.Execute($old2b)
    if ($found2b) {
        $range2b.Text = "sa."
        $range2b.Font.Bold = 1
        $range2b.Font.Bold = 0
    }

    # --- Remove everything from the paragraph right after this one through
    #     the end of the document. This drops the "Brindar...", "Realizar...",
    #     the blank spacer paragraphs, and the "Justificación" section, since
    #     the edited paragraph above is now the last paragraph before the
    #     section properties. ---
    if ($targetParaIndex -lt $d.Paragraphs.Count) {
        $delStart = $d.Paragraphs($targetParaIndex + 1).Range.Start
        $delEnd = $d.Content.End
        $delRange = $d.Range($delStart, $delEnd)
        $delRange.Delete()
    }
}
